# Applies "added ifoCAST full series evaluation" update:
# Columns B:G for rows 2-11 are updated to reflect a shifted evaluation
# window (row 2 drops off, new row 11 data is appended).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$data = @(
    @(-0.2243426211202774, 2.382314769375817, 19.09061844583581, 4.369281227597487, 4.46620300197608, 22),
    @(-0.7120868115500602, 1.933618923903567, 10.88254014902561, 3.298869525917267, 3.300643293046703, 21),
    @(-0.3051420100610381, 1.578682252322016, 10.2251943824045, 3.197685785439917, 3.265784782946779, 20),
    @(-0.2804117665781, 1.683207935177004, 11.04924039551323, 3.324039770446983, 3.402952828034286, 19),
    @(-0.3173356620072107, 1.884213172660337, 10.94775440033029, 3.308739095234057, 3.388969517186058, 18),
    @(-0.2356007239738259, 1.940837229566887, 12.01945225291491, 3.466908169091721, 3.565345857874665, 17),
    @(-0.1538742921933148, 1.930586848953809, 10.98523007623265, 3.314397392623983, 3.41940389010477, 16),
    @(-0.1230455590706185, 1.907399787585887, 12.46868919605704, 3.531103113200893, 3.652819212576486, 15),
    @(-0.09178185018074302, 2.13619364861378, 13.88641155430669, 3.726447578365579, 3.865944277603689, 14),
    @(-0.08394991111291034, 2.334233754498058, 14.90545946209773, 3.860758923074286, 4.017455189220993, 13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
